$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell value: D11 relatedness from "Not Related" to "Related"
$ws.Range("D11").Value = "Related"

# Update the selected cell / active cell on the sheet
$ws.Range("G13").Select()

# Update the workbook window position/size
$excel.Windows.Item(1).Left = 0
$excel.Windows.Item(1).Top = 760
$excel.Windows.Item(1).Width = 34560
$excel.Windows.Item(1).Height = 21580
